$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5393
$ws.Range("F6").Value = 5393
$ws.Range("F7").Value = 186
$ws.Range("F11").Value = 1209
$ws.Range("F12").Value = 6358
$ws.Range("F15").Value = 113
$ws.Range("F16").Value = 3340
$ws.Range("F17").Value = 270
$ws.Range("F18").Value = 108
$ws.Range("F19").Value = 258
$ws.Range("F20").Value = 4063
$ws.Range("F24").Value = 3981
$ws.Range("F25").Value = 193
$ws.Range("F28").Value = 252
$ws.Range("F29").Value = 259
$ws.Range("F32").Value = 136
$ws.Range("F34").Value = 57
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 7110
$ws.Range("F38").Value = 33
$ws.Range("F39").Value = 1173
$ws.Range("F40").Value = 569
$ws.Range("F43").Value = 1431
$ws.Range("F44").Value = 186
$ws.Range("F45").Value = 770
$ws.Range("F46").Value = 3151
$ws.Range("F47").Value = 326
$ws.Range("F49").Value = 796
$ws.Range("F50").Value = 991

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 138
$ws.Range("F25").Value = 837

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 223

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 223
$ws.Range("F8").Value = 5393
$ws.Range("F9").Value = 5393
$ws.Range("F10").Value = 186
$ws.Range("F15").Value = 1209
$ws.Range("F16").Value = 6358
$ws.Range("F19").Value = 113
$ws.Range("F20").Value = 3340
$ws.Range("F21").Value = 270
$ws.Range("F22").Value = 108
$ws.Range("F23").Value = 258
$ws.Range("F24").Value = 4063
$ws.Range("F25").Value = 3981
$ws.Range("F26").Value = 193
$ws.Range("F28").Value = 252
$ws.Range("F29").Value = 259
$ws.Range("F32").Value = 136
$ws.Range("F33").Value = 57
$ws.Range("F36").Value = 7110
$ws.Range("F37").Value = 33
$ws.Range("F38").Value = 1173
$ws.Range("F39").Value = 569
$ws.Range("F43").Value = 1431
$ws.Range("F44").Value = 186
$ws.Range("F45").Value = 770
$ws.Range("F46").Value = 3151
$ws.Range("F47").Value = 326
$ws.Range("F48").Value = 796
$ws.Range("F49").Value = 991
